$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for columns C (sum_removal), D (sum_ha), E (sum_removal_ha)
# Row 7 is intentionally left untouched (no change in the source diff).
$data = @{
    2  = @{ C = 50332051.47490814;  D = 773500144.8783754;  E = 65.07051331299016 }
    3  = @{ C = 781794475.4946911;  D = 4863993091.163775;  E = 160.7310004027239 }
    4  = @{ C = 12417625.785737;    D = 202668507.2253587;  E = 61.27062342216361 }
    5  = @{ C = 240879170.3501668; D = 1667710662.935601;  E = 144.4370271796172 }
    6  = @{ C = 32271138.3268733;  D = 168243183.6822288;  E = 191.8124563538085 }
    8  = @{ C = 501203201.3946835; D = 1573907746.592665;  E = 318.4450946885119 }
    9  = @{ C = 1715879494.422725; D = 4929796264.592522;  E = 348.0629629152745 }
    10 = @{ C = 516414783.6085377; D = 1413441928.104768;  E = 365.3597458375804 }
    11 = @{ C = 894385832.9471077; D = 2315684868.422894;  E = 386.2295103893962 }
    12 = @{ C = 43612244.87203784; D = 179011937.08526;    E = 243.6275791556076 }
    13 = @{ C = 80814206.242;      D = 214025158.8077344;  E = 377.5920863331679 }
    14 = @{ C = 89621261.04626356; D = 1500091513.779432;  E = 59.74386243974256 }
    15 = @{ C = 1078061040.322272; D = 5146429395.837155;  E = 209.477476013621 }
    16 = @{ C = 70465079.67538077; D = 794201029.9429028;  E = 88.7244879050921 }
    17 = @{ C = 485991296.3443565; D = 2437592220.690554;  E = 199.3735015312275 }
    18 = @{ C = 4964423.86935757;  D = 77631574.8399705;   E = 63.94851424296388 }
    19 = @{ C = 57668708.8547;     D = 194082613.4405517;  E = 297.1348532070554 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
}
